# Updates workbook data to reflect the latest scrape snapshot
# (commit: "Update gh-pages to output generated at 456a3b4").
# Bumps "want to go" counts across the 展览 / 演出 / 全部类型 sheets and
# refreshes one event's renamed title + cover image.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (exhibitions) ---
$ws1.Range("F2").Value = 28
$ws1.Range("F5").Value = 86
$ws1.Range("F8").Value = 520
$ws1.Range("F10").Value = 1337
$ws1.Range("F13").Value = 178
$ws1.Range("F14").Value = 20
$ws1.Range("F18").Value = 1690
$ws1.Range("C20").Value = "杭州·热血番&运动番ONLY"
$ws1.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202405/p59Qx5yN1715828421849.jpeg"
$ws1.Range("F21").Value = 265
$ws1.Range("F22").Value = 2953
$ws1.Range("F24").Value = 409
$ws1.Range("F26").Value = 934
$ws1.Range("F27").Value = 1219
$ws1.Range("F29").Value = 2855
$ws1.Range("F30").Value = 1651
$ws1.Range("F33").Value = 679
$ws1.Range("F35").Value = 1909
$ws1.Range("F36").Value = 904
$ws1.Range("F37").Value = 1919
$ws1.Range("F38").Value = 209
$ws1.Range("F39").Value = 43
$ws1.Range("F40").Value = 58
$ws1.Range("F43").Value = 903
$ws1.Range("F44").Value = 812
$ws1.Range("F45").Value = 1051
$ws1.Range("F46").Value = 133
$ws1.Range("F47").Value = 450
$ws1.Range("F48").Value = 232
$ws1.Range("F49").Value = 3373

# --- Sheet "演出" (performances) ---
$ws2.Range("F7").Value = 14
$ws2.Range("F12").Value = 811

# --- Sheet "全部类型" (all types, combined view) ---
$ws4.Range("F4").Value = 86
$ws4.Range("F9").Value = 520
$ws4.Range("F11").Value = 1338
$ws4.Range("F14").Value = 178
$ws4.Range("F18").Value = 1690
$ws4.Range("C20").Value = "杭州·热血番&运动番ONLY"
$ws4.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202405/p59Qx5yN1715828421849.jpeg"
$ws4.Range("F21").Value = 265
$ws4.Range("F22").Value = 2953
$ws4.Range("F24").Value = 409
$ws4.Range("F25").Value = 14
$ws4.Range("F26").Value = 1219
$ws4.Range("F27").Value = 2855
$ws4.Range("F28").Value = 1651
$ws4.Range("F31").Value = 811
$ws4.Range("F34").Value = 1909
$ws4.Range("F36").Value = 904
$ws4.Range("F37").Value = 1919
$ws4.Range("F38").Value = 43
$ws4.Range("F39").Value = 58
$ws4.Range("F41").Value = 903
$ws4.Range("F42").Value = 812
$ws4.Range("F43").Value = 1051
$ws4.Range("F44").Value = 133
$ws4.Range("F45").Value = 450
$ws4.Range("F47").Value = 232
$ws4.Range("F48").Value = 3373
